# Updates cryptos list values (Price / Volume(1h) columns) and two row swaps
# (PEPE/Polkadot rows 22-23, Mantle -> Filecoin row 51) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.787.11"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "3.343.76"
$ws.Range("E3").Value = "  +8.78%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "257.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "3.342.57"
$ws.Range("E10").Value = "  +8.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "97.482.91"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000246"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.42%  "
$ws.Range("D16").Value = "3.962.61"
$ws.Range("E16").Value = "  +8.44%  "
$ws.Range("E17").Value = "  +5.38%  "
$ws.Range("D18").Value = "3.342.31"
$ws.Range("E18").Value = "  +8.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "485.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.25%  "
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000207"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("E25").Value = "  +2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("D28").Value = "3.585.21"
$ws.Range("E28").Value = "  +10.25%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.79%  "
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.08%  "
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.448"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.11%  "
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.82%  "
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("E45").Value = "  +17.88%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("E49").Value = "  +8.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.44%  "
